$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A56").Value = "5EJY4Q"
$ws.Range("B56").Value = "Kit de engranaje de fusor Ricoh"
$ws.Range("C56").Value = "MP C2003 C2004 C2503 C3003 C3503 C4503 C5503 C6003"
$ws.Range("D56").Value = 63000
$ws.Range("E56").Value = 150000
$ws.Range("F56").Value = 2
$ws.Range("G56").Value = 1
$ws.Range("H56").Formula = "=(E56-D56)*G56"
$ws.Range("I56").Formula = "=D56*F56"
$ws.Range("J56").Value = 126000
